# Add an "antibodies_path" column to the IMC metadata sheet, between
# data_precision_bytes (AL) and the existing contributors_path (was AM,
# now shifts to AN) / data_path (was AN, now shifts to AO) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at AM; this pushes the old AM (contributors_path)
# to AN and the old AN (data_path) to AO. Values move with the insert, but
# cell comments stay pinned to their original column, so they need to be
# re-applied below.
$ws.Range("AM1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("AM1").Value = "antibodies_path"

# Re-point the header-row comments so they travel with their header text.
$ws.Range("AM1").Comment.Text("Relative path to file with antibody information for this dataset.")
$ws.Range("AN1").Comment.Text("Relative path to file with ORCID IDs for contributors for this dataset.")
$ws.Range("AO1").AddComment("Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.")

# Match the bold/centered/wrapped header style used by every other column 1
# header cell (the new AO1 starts out with the default, unstyled format).
$ws.Range("AO1").Font.Bold = $true
$ws.Range("AO1").HorizontalAlignment = -4108
$ws.Range("AO1").WrapText = $true
